$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: Monitor / Monitoring computer program running process ---
$ws.Cells.Item(3, 2).Value() = "Monitoring computer program running process"
$ws.Cells.Item(3, 4).ClearContents()
$ws.Cells.Item(3, 5).Value() = 2

# --- Row 5: Identify which application processes are games ---
$ws.Cells.Item(5, 2).Value() = "Identify which application processes are games"

# --- Row 6: unchanged (Draw statistics graph of game time) ---

# --- Row 7: unchanged (Allow parents to set daily or weekly limits for gaming time) ---

# --- Row 8: Automatically interrupt the game process when reach game time limit ---
$ws.Cells.Item(8, 2).Value() = "Automatically interrupt the game process when reach game time limit"
$ws.Cells.Item(8, 4).Value() = 1.5

# --- Row 9: Remind children when approaching the time limit, moved estimate from E to D ---
$ws.Cells.Item(9, 5).ClearContents()
$ws.Cells.Item(9, 4).Value() = 0.5

# --- Row 10: Parents can remotely interrupt the game process through mobile phones ---
$ws.Cells.Item(10, 2).Value() = "Parents can remotely interrupt the game process through mobile phones"

# --- Row 4: Record the running time of each game application process ---
$ws.Cells.Item(4, 2).Value() = "Record the running time of each game application process"
$ws.Cells.Item(4, 4).ClearContents()
$ws.Cells.Item(4, 5).Value() = 1

# --- Row 11: new task row ---
$ws.Cells.Item(11, 1).Value() = "Cloud(Telegram)"
$ws.Cells.Item(11, 2).Value() = "Transfer commands and data through the Cloud (e.g. Telegram)"
$ws.Cells.Item(11, 3).Value() = "H"
$ws.Cells.Item(11, 4).Value() = 0.5
$ws.Cells.Item(11, 5).Value() = 1.5

# --- Column widths (best-fit for the now-wider text; A=15.2, B=65 chars) ---
$ws.Columns.Item(1).ColumnWidth = 14.428571428571429
$ws.Columns.Item(2).ColumnWidth = 64.28571428571429

# --- Update the selected cell to match the saved view state ---
$ws.Range("F9").Select()
